$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("predictability")
$co = $ws.ChartObjects().Add(100, 100, 300, 200)
$ch = $co.Chart
$ch.ChartType = 51
$ch.SetSourceData($ws.Range("B61:D64"))
Write-Host "done"
